$wb = $excel.ActiveWorkbook

# Update "展览" sheet (F2: 105 -> 109, F3: 21 -> 23)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 109
$ws1.Range("F3").Value = 23

# Update "全部类型" sheet (F2: 105 -> 109, F3: 21 -> 23)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 109
$ws4.Range("F3").Value = 23
